$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; all existing columns (and their
# formulas) shift one to the right automatically.
$ws.Range("A1").EntireColumn.Insert()

# New first column: snowdepth measurements for every row.
$ws.Range("A1").Value = "snowdepth"
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 17
$ws.Range("A4").Value = 32
$ws.Range("A5").Value = 45
$ws.Range("A6").Value = 48
$ws.Range("A7").Value = 68
$ws.Range("A8").Value = 74
$ws.Range("A9").Value = 84
$ws.Range("A10").Value = 89
$ws.Range("A11").Value = 118

# Rename the shifted headers (old "Schneehöhe" / "Mittelwert" -> new text).
$ws.Range("B1").Value = "snowheighth"
$ws.Range("H1").Value = "mean"

$ws.Range("H1").Select() | Out-Null
